# Year to Year-Group Linkage.xlsx edit script
# 1) Duplicate the current "main" sheet (which holds the OLD 2000-2019 data) to
#    become the "old" tab, fix its one data bug (E11: 1009 -> 2009), and adjust
#    its view state.
# 2) Rewrite the original "main" sheet in place with the refreshed 2000-2020
#    data set (new multi-year groupings), adding row 22.
# 3) Tidy up sheet order / selection / the now-unused blank "Sheet1".

$wb = $excel.ActiveWorkbook

$mainOld = $wb.Worksheets.Item("main")

# --- Step 1: snapshot of the current (old) "main" sheet becomes "old" ---
$mainOld.Copy($null, $mainOld)
$oldWs = $wb.Worksheets.Item("main (2)")
$oldWs.Name = "old"

# Fix the data bug on the "old" tab.
$oldWs.Range("E11").Value = 2009

# "old" is no longer the active/selected tab.
$oldWs.Range("B7").Select()

# --- Step 2: remove the now-unused blank placeholder sheet ---
$wb.Worksheets.Item("Sheet1").Delete()

# --- Step 3: rewrite "main" with the refreshed data set ---
$ws = $mainOld

# Extend formatting down to the new row 22 before writing values, by copying
# the formatting of row 21.
$ws.Range("A21:E21").Copy()
$ws.Range("A22:E22").PasteSpecial(-4122)

$ws.Range("A2").Value = 2000
$ws.Range("B2").Value = "'2000"
$ws.Range("C2").ClearContents()
$ws.Range("D2").Value = "2000-2002"
$ws.Range("E2").Value = 2001

$ws.Range("A3").Formula = "=A2+1"
$ws.Range("B3").Value = "2001-2005"
$ws.Range("C3").Value = 2003
$ws.Range("D3").Value = "2000-2002"
$ws.Range("E3").Value = 2001

$ws.Range("A4").Formula = "=A3+1"
$ws.Range("B4").Value = "2001-2005"
$ws.Range("C4").Value = 2003
$ws.Range("D4").Value = "2000-2002"
$ws.Range("E4").Value = 2001

$ws.Range("A5").Formula = "=A4+1"
$ws.Range("B5").Value = "2001-2005"
$ws.Range("C5").Value = 2003
$ws.Range("D5").Value = "2003-2005"
$ws.Range("E5").Value = 2004

$ws.Range("A6").Formula = "=A5+1"
$ws.Range("B6").Value = "2001-2005"
$ws.Range("C6").Value = 2003
$ws.Range("D6").Value = "2003-2005"
$ws.Range("E6").Value = 2004

$ws.Range("A7").Formula = "=A6+1"
$ws.Range("B7").Value = "2001-2005"
$ws.Range("C7").Value = 2003
$ws.Range("D7").Value = "2003-2005"
$ws.Range("E7").Value = 2004

$ws.Range("A8").Formula = "=A7+1"
$ws.Range("B8").Value = "2006-2010"
$ws.Range("C8").Value = 2008
$ws.Range("D8").Value = "2006-2008"
$ws.Range("E8").Value = 2007

$ws.Range("A9").Formula = "=A8+1"
$ws.Range("B9").Value = "2006-2010"
$ws.Range("C9").Value = 2008
$ws.Range("D9").Value = "2006-2008"
$ws.Range("E9").Value = 2007

$ws.Range("A10").Formula = "=A9+1"
$ws.Range("B10").Value = "2006-2010"
$ws.Range("C10").Value = 2008
$ws.Range("D10").Value = "2006-2008"
$ws.Range("E10").Value = 2007

$ws.Range("A11").Formula = "=A10+1"
$ws.Range("B11").Value = "2006-2010"
$ws.Range("C11").Value = 2008
$ws.Range("D11").Value = "2009-2011"
$ws.Range("E11").Value = 2010

$ws.Range("A12").Formula = "=A11+1"
$ws.Range("B12").Value = "2006-2010"
$ws.Range("C12").Value = 2008
$ws.Range("D12").Value = "2009-2011"
$ws.Range("E12").Value = 2010

$ws.Range("A13").Formula = "=A12+1"
$ws.Range("B13").Value = "2011-2015"
$ws.Range("C13").Value = 2013
$ws.Range("D13").Value = "2009-2011"
$ws.Range("E13").Value = 2010

$ws.Range("A14").Formula = "=A13+1"
$ws.Range("B14").Value = "2011-2015"
$ws.Range("C14").Value = 2013
$ws.Range("D14").Value = "2012-2014"
$ws.Range("E14").Value = 2013

$ws.Range("A15").Formula = "=A14+1"
$ws.Range("B15").Value = "2011-2015"
$ws.Range("C15").Value = 2013
$ws.Range("D15").Value = "2012-2014"
$ws.Range("E15").Value = 2013

$ws.Range("A16").Formula = "=A15+1"
$ws.Range("B16").Value = "2011-2015"
$ws.Range("C16").Value = 2013
$ws.Range("D16").Value = "2012-2014"
$ws.Range("E16").Value = 2013

$ws.Range("A17").Formula = "=A16+1"
$ws.Range("B17").Value = "2011-2015"
$ws.Range("C17").Value = 2013
$ws.Range("D17").Value = "2015-2017"
$ws.Range("E17").Value = 2016

$ws.Range("A18").Formula = "=A17+1"
$ws.Range("B18").Value = "2016-2020"
$ws.Range("C18").Value = 2018
$ws.Range("D18").Value = "2015-2017"
$ws.Range("E18").Value = 2016

$ws.Range("A19").Value = 2017
$ws.Range("B19").Value = "2016-2020"
$ws.Range("C19").Value = 2018
$ws.Range("D19").Value = "2015-2017"
$ws.Range("E19").Value = 2016

$ws.Range("A20").Value = 2018
$ws.Range("B20").Value = "2016-2020"
$ws.Range("C20").Value = 2018
$ws.Range("D20").Value = "2018-2020"
$ws.Range("E20").Value = 2019

$ws.Range("A21").Value = 2019
$ws.Range("B21").Value = "2016-2020"
$ws.Range("C21").Value = 2018
$ws.Range("D21").Value = "2018-2020"
$ws.Range("E21").Value = 2019

$ws.Range("A22").Value = 2020
$ws.Range("B22").Value = "2016-2020"
$ws.Range("C22").Value = 2018
$ws.Range("D22").Value = "2018-2020"
$ws.Range("E22").Value = 2019

# Column widths on "main": columns B:C width 14, column D width 18.28515625,
# column A back to the workbook default (no explicit style/width override).
$ws.Columns("B:C").ColumnWidth = 14
$ws.Columns("D:D").ColumnWidth = 18.28515625

# "main" stays the active / selected tab, cursor parked just past the data.
$ws.Range("A23").Select()

# Make sure "main" is first, "old" second, and "main" is the active sheet.
$ws.Move($wb.Worksheets.Item(1))
$ws.Activate()
